$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.227.28"
$ws.Range("E2").Value = "  +0.21%  "

$ws.Range("D3").Value = "1.833.36"
$ws.Range("E3").Value = "  -0.48%  "

$ws.Range("D4").Value = "0.9985"
$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").Value = "242.59"
$ws.Range("E5").Value = "  -0.75%  "

$ws.Range("D6").Value = "0.6239"
$ws.Range("E6").Value = "  -0.38%  "

$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  -0.24%  "

$ws.Range("D8").Value = "0.07380"
$ws.Range("E8").Value = "  -1.72%  "

$ws.Range("D9").Value = "0.2904"
$ws.Range("E9").Value = "  -1.22%  "

$ws.Range("D10").Value = "23.19"
$ws.Range("E10").Value = "  -0.56%  "

$ws.Range("D11").Value = "0.07675"
$ws.Range("E11").Value = "  -0.62%  "

$ws.Range("D12").Value = "1.824.42"
$ws.Range("E12").Value = "  -1.80%  "

$ws.Range("D14").Value = "0.6679"
$ws.Range("E14").Value = "  -1.31%  "

$ws.Range("D15").Value = "82.63"
$ws.Range("E15").Value = "  -0.59%  "

$ws.Range("D16").Value = "0.000008970"
$ws.Range("E16").Value = "  -3.26%  "

$ws.Range("D17").Value = "5.880"
$ws.Range("E17").Value = "  -1.63%  "

$ws.Range("D18").Value = "29.193.02"
$ws.Range("E18").Value = "  +0.10%  "

$ws.Range("D19").Value = "2.073.36"
$ws.Range("E19").Value = "  -2.74%  "

$ws.Range("D20").Value = "236.24"

$ws.Range("D21").Value = "12.48"
$ws.Range("E21").Value = "  -1.79%  "

$ws.Range("D22").Value = "0.9996"
$ws.Range("E22").Value = "  -0.31%  "

$ws.Range("D23").Value = "7.392"
$ws.Range("E23").Value = "  +2.72%  "

$ws.Range("D24").Value = "0.9998"
$ws.Range("E24").Value = "  -0.26%  "

$ws.Range("D25").Value = "158.27"
$ws.Range("E25").Value = "  -1.48%  "

$ws.Range("D26").Value = "0.1410"
$ws.Range("E26").Value = "  +1.25%  "

$ws.Range("D27").Value = "8.540"
$ws.Range("E27").Value = "  -0.35%  "

$ws.Range("D28").Value = "17.67"
$ws.Range("E28").Value = "  -1.46%  "

$ws.Range("D29").Value = "1.483"
$ws.Range("E29").Value = "  -1.27%  "

$ws.Range("D30").Value = "0.05813"
$ws.Range("E30").Value = "  +4.48%  "

$ws.Range("D31").Value = "4.101"
$ws.Range("E31").Value = "  -1.23%  "

$ws.Range("D32").Value = "4.086"
$ws.Range("E32").Value = "  -2.57%  "

$ws.Range("D33").Value = "1.204"
$ws.Range("E33").Value = "  -0.21%  "

$ws.Range("D34").Value = "1.871"
$ws.Range("E34").Value = "  +0.72%  "

$ws.Range("E35").Value = "  -2.50%  "

$ws.Range("D36").Value = "1.141"
$ws.Range("E36").Value = "  -0.52%  "

$ws.Range("D37").Value = "2.609"

$ws.Range("D38").Value = "2.848"
$ws.Range("E38").Value = "  +2.62%  "

$ws.Range("D39").Value = "1.225.61"
$ws.Range("E39").Value = "  -0.31%  "

$ws.Range("D40").Value = "0.01758"
$ws.Range("E40").Value = "  -1.89%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "6.287"
$ws.Range("E41").Value = "  -4.40%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "0.9152"
$ws.Range("E42").Value = "  +1.31%  "

$ws.Range("E43").Value = "  -0.13%  "

$ws.Range("D44").Value = "101.86"
$ws.Range("E44").Value = "  -0.44%  "

$ws.Range("D45").Value = "1.977.39"
$ws.Range("E45").Value = "  -2.15%  "

$ws.Range("D46").Value = "65.13"
$ws.Range("E46").Value = "  -1.91%  "

$ws.Range("D47").Value = "0.5038"
$ws.Range("E47").Value = "  -1.24%  "

$ws.Range("E48").Value = "  -3.66%  "

$ws.Range("D49").Value = "0.4031"
$ws.Range("E49").Value = "  -1.46%  "

$ws.Range("D50").Value = "9.098"
$ws.Range("E50").Value = "  -0.44%  "

$ws.Range("D51").Value = "0.1133"
$ws.Range("E51").Value = "  +2.71%  "
